$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '49.853.63'
$ws.Range("E2").Value = '  +3.73%  '

$ws.Range("D3").Value = '2.642.83'
$ws.Range("E3").Value = '  +5.74%  '

$ws.Range("E4").Value = '  +0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '113.76'
$ws.Range("E5").Value = '  +7.47%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '326.15'
$ws.Range("E6").Value = '  +2.03%  '

$ws.Range("E7").Value = '  +1.29%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'

$ws.Range("E9").Value = '  +2.59%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.91'
$ws.Range("E10").Value = '  +5.63%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.15'
$ws.Range("E11").Value = '  +0.64%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0820'

$ws.Range("E13").Value = '  +0.92%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.34'
$ws.Range("E14").Value = '  +3.44%  '

$ws.Range("D15").Value = '3.054.51'
$ws.Range("E15").Value = '  +5.65%  '

$ws.Range("D16").Value = '2.666.57'
$ws.Range("E16").Value = '  +6.52%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.868'
$ws.Range("E17").Value = '  +4.21%  '

$ws.Range("D18").Value = '49.759.74'
$ws.Range("E18").Value = '  +3.92%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.14'
$ws.Range("E19").Value = '  +0.75%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.77'
$ws.Range("E20").Value = '  +1.70%  '

$ws.Range("E21").Value = '  -0.48%  '

$ws.Range("D22").Value = '0.0₃0956'
$ws.Range("E22").Value = '  +2.41%  '

$ws.Range("B23").Value = 'BitcoinCash'
$ws.Range("C23").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '277.78'
$ws.Range("E23").Value = '  +1.79%  '

$ws.Range("B24").Value = 'Litecoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '72.05'
$ws.Range("E24").Value = '  +1.17%  '

$ws.Range("E25").Value = '  +2.21%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.72'
$ws.Range("E26").Value = '  +3.65%  '

$ws.Range("B28").Value = 'Cosmos'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.98'
$ws.Range("E28").Value = '  +2.49%  '

$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.23'
$ws.Range("E29").Value = '  -2.94%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '35.97'
$ws.Range("E30").Value = '  +3.22%  '

$ws.Range("E31").Value = '  -1.44%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '50.33'
$ws.Range("E32").Value = '  +2.13%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.42'
$ws.Range("E33").Value = '  +2.32%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '19.44'
$ws.Range("E34").Value = '  +1.52%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0805'
$ws.Range("E35").Value = '  +3.92%  '

$ws.Range("E36").Value = '  -0.05%  '

$ws.Range("E37").Value = '  +6.32%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.81'
$ws.Range("E38").Value = '  +5.25%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.07'
$ws.Range("E39").Value = '  +7.07%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '125.49'
$ws.Range("E40").Value = '  +3.19%  '

$ws.Range("E41").Value = '  +1.56%  '

$ws.Range("E42").Value = '  +1.40%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '22.03'
$ws.Range("E43").Value = '  -1.09%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0314'
$ws.Range("E44").Value = '  +3.19%  '

$ws.Range("D45").Value = '2.071.09'
$ws.Range("E45").Value = '  +3.31%  '

$ws.Range("E46").Value = '  +4.86%  '

$ws.Range("E47").Value = '  +14.85%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.97'
$ws.Range("E48").Value = '  +4.80%  '

$ws.Range("E49").Value = '  +2.04%  '

$ws.Range("E50").Value = '  +3.70%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '59.53'
$ws.Range("E51").Value = '  +5.75%  '
